# Insert a new price-record row at row 84 (pushing the existing rows 84-184
# down to 85-185) and populate it with the new Alcachofa ("Española",
# Primera) observation dated 2022-07-21 (serial 44763).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(84).Insert()

$ws.Cells.Item(84, 1).Value = 10
$ws.Cells.Item(84, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(84, 3).Value = 'La Araucanía'
$ws.Cells.Item(84, 4).Value = 44763
$ws.Cells.Item(84, 5).Value = 9
$ws.Cells.Item(84, 6).Value = 100112013
$ws.Cells.Item(84, 7).Value = 'Alcachofa'
$ws.Cells.Item(84, 8).Value = 'Española'
$ws.Cells.Item(84, 9).Value = 'Primera'
$ws.Cells.Item(84, 10).Value = 115
$ws.Cells.Item(84, 11).Value = 18000
$ws.Cells.Item(84, 12).Value = 20000
$ws.Cells.Item(84, 13).Value = 19130
$ws.Cells.Item(84, 14).Value = '$/caja 30 unidades'
$ws.Cells.Item(84, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(84, 16).Value = 638
$ws.Cells.Item(84, 17).Value = 30
$ws.Cells.Item(84, 18).Value = 'Hortaliza'
